# Refresh cryptos list: updated prices / 1h volume %, and
# re-sorted three coin-pairs (rows 11/12, 20/21, 47/48/49) to match the
# latest coinranking.com snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value into a cell without Excel
# re-interpreting numeric-looking strings (e.g. "6.210", "1.108.91")
# as numbers/dates - mirrors how these cells were originally stored
# (inline text), and leaves the cell style untouched afterwards.
function Set-TextValue($addr, $val) {
  $c = $ws.Range($addr)
  $c.NumberFormat = "@"
  $c.Value = $val
  $c.ClearFormats()
}

Set-TextValue "D2" "29.955.32"

Set-TextValue "D3" "1.893.60"

Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue "D5" "0.7784"
$ws.Range("E5").Value = "  +0.42%  "

Set-TextValue "D6" "243.92"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  -0.01%  "

Set-TextValue "D8" "0.3136"
$ws.Range("E8").Value = "  +0.47%  "

Set-TextValue "D9" "25.85"
$ws.Range("E9").Value = "  +2.53%  "

Set-TextValue "D10" "0.07275"
$ws.Range("E10").Value = "  +1.65%  "

$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D11" "2.164.32"
$ws.Range("E11").Value = "  +15.96%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D12" "0.08684"
$ws.Range("E12").Value = "  +7.83%  "

Set-TextValue "D13" "0.7742"
$ws.Range("E13").Value = "  +1.47%  "

$ws.Range("E14").Value = "  -0.59%  "

Set-TextValue "D15" "94.53"
$ws.Range("E15").Value = "  +2.59%  "

Set-TextValue "D16" "6.210"
$ws.Range("E16").Value = "  +1.15%  "

Set-TextValue "D17" "30.119.58"
$ws.Range("E17").Value = "  +1.12%  "

$ws.Range("E18").Value = "  -0.02%  "

Set-TextValue "D19" "245.94"
$ws.Range("E19").Value = "  +1.15%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D20" "2.288.77"
$ws.Range("E20").Value = "  +9.02%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D21" "0.000007880"
$ws.Range("E21").Value = "  +1.70%  "

Set-TextValue "D22" "8.203"
$ws.Range("E22").Value = "  +1.63%  "

$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("E24").Value = "  +0.03%  "

Set-TextValue "D25" "0.1663"
$ws.Range("E25").Value = "  +3.76%  "

Set-TextValue "D26" "9.503"
$ws.Range("E26").Value = "  +1.28%  "

Set-TextValue "D27" "163.39"
$ws.Range("E27").Value = "  +1.29%  "

Set-TextValue "D28" "18.86"
$ws.Range("E28").Value = "  +0.83%  "

Set-TextValue "D29" "2.053"
$ws.Range("E29").Value = "  +0.44%  "

Set-TextValue "D30" "1.433"
$ws.Range("E30").Value = "  +0.85%  "

Set-TextValue "D31" "1.543"
$ws.Range("E31").Value = "  -0.29%  "

Set-TextValue "D32" "4.512"
$ws.Range("E32").Value = "  +0.99%  "

Set-TextValue "D33" "4.131"
$ws.Range("E33").Value = "  +0.84%  "

Set-TextValue "D34" "0.05494"
$ws.Range("E34").Value = "  -0.53%  "

Set-TextValue "D35" "1.249"
$ws.Range("E35").Value = "  -0.90%  "

Set-TextValue "D36" "0.7548"
$ws.Range("E36").Value = "  +1.47%  "

Set-TextValue "D37" "1.001"
$ws.Range("E37").Value = "  +0.45%  "

Set-TextValue "D38" "2.685"
$ws.Range("E38").Value = "  +2.50%  "

Set-TextValue "D39" "0.01961"
$ws.Range("E39").Value = "  +2.55%  "

Set-TextValue "D40" "2.789"
$ws.Range("E40").Value = "  +0.28%  "

Set-TextValue "D41" "0.4515"
$ws.Range("E41").Value = "  +2.30%  "

Set-TextValue "D42" "1.108.91"
$ws.Range("E42").Value = "  -2.87%  "

Set-TextValue "D43" "73.71"
$ws.Range("E43").Value = "  +0.33%  "

Set-TextValue "D44" "6.079"
$ws.Range("E44").Value = "  +3.95%  "

Set-TextValue "D45" "0.8531"
$ws.Range("E45").Value = "  +0.30%  "

Set-TextValue "D46" "1.000"

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D47" "103.23"
$ws.Range("E47").Value = "  -0.46%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D48" "1.882"
$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D49" "2.158.13"
$ws.Range("E49").Value = "  +7.60%  "

Set-TextValue "D50" "7.605"
$ws.Range("E50").Value = "  +2.35%  "

Set-TextValue "D51" "9.868"
$ws.Range("E51").Value = "  -1.04%  "
